$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorganize the account-statement rows (16-19): previously grouped by period
# (2105 then 2107) for each worker; now grouped by worker (ANDRES then FRANKLIN)
# with each worker's two periods (2107 then 2105) listed together.

# Row 16: ANDRES GILBERTO VEGA MARRUGO - period 2107
$ws.Range("C16").Value = "73119097"
$ws.Range("D16").Value = "ANDRES GILBERTO VEGA MARRUGO"
$ws.Range("E16").Value = "2107"
$ws.Range("F16").Value = 36341

# Row 17: ANDRES GILBERTO VEGA MARRUGO - period 2105
$ws.Range("C17").Value = "73119097"
$ws.Range("D17").Value = "ANDRES GILBERTO VEGA MARRUGO"
$ws.Range("E17").Value = "2105"
$ws.Range("F17").Value = 36000

# Row 18: FRANKLIN ENRIQUE VEGA MARRUGO - period 2107
$ws.Range("C18").Value = "73070176"
$ws.Range("D18").Value = "FRANKLIN ENRIQUE VEGA MARRUGO"
$ws.Range("E18").Value = "2107"
$ws.Range("F18").Value = 36341

# Row 19: FRANKLIN ENRIQUE VEGA MARRUGO - period 2105
$ws.Range("C19").Value = "73070176"
$ws.Range("D19").Value = "FRANKLIN ENRIQUE VEGA MARRUGO"
$ws.Range("E19").Value = "2105"
$ws.Range("F19").Value = 35112

$wb.Save()
